$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...muon nhan manh su khao khat duoc thuong trang..." ->
#           "...muon nhan manh " | "niem" | " khao khat duoc thuong trang..."
# i.e. replace the word "su" ("sự") with "niem" ("niềm") while keeping the
# surrounding text as separate runs (matching the target XML diff, which
# splits the single original run into three runs).
# ---------------------------------------------------------------------------

# Locate "sự" inside "sự khao khát" (unique in the document).
$suRange = $d.Content
$suRange.Find.Execute("sự khao khát", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$suStart = $suRange.Start
$suEnd = $suStart + 2

# Locate the end of the original run ("... bản thân mình").
$tailRange = $d.Content
$tailRange.Find.Execute("bản thân mình", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$runEnd = $tailRange.End

# Replace "sự" with "niềm". Word naturally turns the replaced word into its
# own run, split away from the text before it.
$repl = $d.Range($suStart, $suEnd)
$repl.Text = "niềm"
# Re-toggling a character property on the new run forces it to stay a
# separate run (rather than being silently re-merged with its neighbours).
$repl.Bold = $true
$repl.Bold = $false

# "niềm" grew the text by 2 characters (2 -> 4), so the tail portion now
# starts right after the inserted word and the original run now ends 2
# characters later than before.
$tailStart = $repl.End
$newRunEnd = $runEnd + ("niềm".Length - 2)

$tailPart = $d.Range($tailStart, $newRunEnd)
$tailPart.Bold = $true
$tailPart.Bold = $false

# ---------------------------------------------------------------------------
# Change 2: collapse the run of "," " " "tình yêu thiên nhiên," " "
#           "phong thái ung dung," " " "lạc quan của bác thật cao cả làm
#           sao!(13)" runs (directly following the underlined "Ôi" run)
#           into a single run with the same concatenated text.
# ---------------------------------------------------------------------------
$mergeText = ", tình yêu thiên nhiên, phong thái ung dung, lạc quan của bác thật cao cả làm sao!(13)"
$d.Content.Find.Execute($mergeText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $mergeText, 2) | Out-Null
